$d = $word.ActiveDocument

# Fix Semi's spelling error: "Nearly ever Blue Bike" -> "Nearly every Blue Bike"
$d.Content.Find.Execute("Nearly ever Blue Bike", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nearly every Blue Bike", 2)
